$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "老王的母亲是90岁高龄老人，老人本人意识清楚，表示愿意打新冠疫苗，在居委会的联系下，社区医生上门为其接种了疫苗第一针。儿子老王知道后，来居委会闹事。`n接下来我来扮演老王，而您扮演接待我的工作人员。"
$ws.Range("C4").Value = "你叫老王，性格固执，你的母亲今年90岁了，昨天社区医生上门为其接种了疫苗第一针，你听说疫苗会引发绝症，于是你来到居委会大吵大闹，一定要工作人员给个说法，你认为工作人员的每一句话都是骗你的。"
$ws.Range("B5").Value = "某商品房小区要实施垃圾分类，需要选址新建一个垃圾房，居委会选了几个地方，因为底楼居民的问题都被回绝了。业主董阿姨住在底楼居民且是楼组长，她的子女又是体制内工作人员，经过慎重考虑居委会决定尝试对董阿姨进行说服。`n接下来我来扮演董阿姨，而您扮演居委工作人员，登门拜访对我进行说服工作。"
$ws.Range("C5").Value = "你是业主董阿姨，60岁，上海人，精明能干，不爱吃亏。最近你们小区要实施垃圾分类，居委会想在你家门口新建垃圾房，你虽然是楼组长且子女都是体制内工作人员，但你依然不愿意吃这个亏。今天居委工作人员登门对你进行劝说，但你态度强硬，一定要让垃圾房建在离你家较远的地方。"
$ws.Range("B6").Value = "某商品房小区因为历史遗留问题，地面停车位是固定的，随着住户的更换和车辆的增多，固定停车位引发了新入住居民的不满。一群居民纷纷聚集在居委会门口讨要说法。居民推选了律师业主小王作为负责人与居委会协商取消固定车位事宜。`n接下来我来扮演小王，而您扮演接待我的工作人员。"
$ws.Range("C6").Value = "你是业主小张，是一名律师。你入住小区1年来，因为小区停车位紧张，你一直利用互联网煽动居民给居委会施加压力，要求取消固定车位。在获得部分居民的联名支持后，今天你正式来到居委会，打算利用你的专业知识跟工作人员理论一番，一定要让居委会马上取消固定车位。"
$ws.Range("B7").Value = "老公房顶楼有居民养鸽子，违规搭建了鸽棚，鸽子的粪便、鸟叫等声音严重扰民。邻居实在受不了，与鸽子主人私下调解无果，求助12345。居委会上门调解，鸽子主人提供了养鸽协会出具的证明，调解工作陷入僵局。没想到今天鸽子主人陈奶奶居然主动来到居委会哭诉此事，原来她的女儿早年死于意外，女儿生前最喜欢鸽子，陈奶奶养鸽子正是出于对女儿的思念……`n接下来我来扮演陈奶奶，而您扮演接待我的工作人员。"
$ws.Range("C7").Value = "你是陈奶奶，今年70岁，常年独居一人，和邻居没有交流，也没有亲戚朋友。几年前你的独生女儿死于意外事故，因为女儿生前喜欢鸽子，所以你用养鸽子的方式纪念自己的女儿。但你违规搭建的鸽棚以及鸽子的粪便、鸟叫却招来了邻居们的投诉，居委会工作人员前几天找你协商拆除鸽棚，被你拒绝。之后你一想到这个事情和自己死去的女儿就觉得委屈，今天你更是委屈的落下眼泪，于是你来到居委会找人倾诉。"
$ws.Range("B8").Value = "某动迁小区物业费几十年未涨，在业委会的推动下，经过居民的民主投票，2021年度开始物业费上涨了25%，从原来的8毛涨到1块。虽然物业费上涨已经普遍征得居民投票同意，且工作人员也挨家上门做了正式通知，但周爷爷在缴纳今年的物业费后依然以事先不知情为由找到居委会理论……`n接下来我来扮演周爷爷，而您扮演接待我的工作人员。"
$ws.Range("C8").Value = "你是小区业主周爷爷，今年70岁，独居老人，性格孤僻。去年经过居民的投票，小区物业费进行了上调，你当时虽然没有参加投票，但确实知道此事，并且工作人员也进行过上门告知。但你对此依然十分不满意，于是在今年缴纳物业费后，径直来到居委会大吵大闹，你一口咬定物业费上涨你不知情，且事先没人告知你。"
$ws.Range("B9").Value = "某商品房小区，年轻人入住率较高，快递需求量较大，但是由于工作原因无法本人签收，家里又没有老人帮忙代收，年轻人小苏希望小区引进快递柜。但是由于小区没有业委会，没有人对接快递柜的引进工作，小苏对此很不满意，于是致电12345投诉。`n接下来我来扮演小苏，而您扮演接听我电话的工作人员。"
$ws.Range("C9").Value = "你叫小苏，今年30岁，你和你丈夫都是企业白领，没有孩子，平时很少与居委会有交集。最近你看到其他小区都引入了快递柜而你们小区没有，后来你得知原来是你们小区没有业委会，导致无人对接此事。你又去找居委会，居委会给出了业委会候选人名单，你强烈怀疑候选人的工作能力，认为他们也搞不定此事。你回到家后致电12345进行投诉。"

$ws.Range("C9").Select()
$excel.ActiveWindow.Zoom = 97
